# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# This script updates the response/error sheets of the workbook so that the
# detailed field-by-field documentation rows are collapsed into a single
# "schema" reference row pointing at the relevant named schema.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: turn row 3 of a sheet (which currently documents the "dateTime"
# field) into a single schema-reference row, and delete any further detail
# rows that used to follow it (commandRef / errorCode / etc).
# ---------------------------------------------------------------------------
function Set-SchemaRow($ws, $lastDetailRow, $schemaName) {
    # Column B: field name -> schema name
    $ws.Cells.Item(3, 2).Value = $schemaName
    # Column D: description -> cleared
    $ws.Cells.Item(3, 4).Value = ""
    # Column E: type string -> schema
    $ws.Cells.Item(3, 5).Value = "schema"
    # Column G: schema name column populated
    $ws.Cells.Item(3, 7).Value = $schemaName
    # Column I: mandatory M -> Yes
    $ws.Cells.Item(3, 9).Value = "Yes"
    # Column L: format pattern -> cleared
    $ws.Cells.Item(3, 12).Value = ""
    # Column O: example -> cleared
    $ws.Cells.Item(3, 15).Value = ""

    # Remove any extra detail rows (rows 4..lastDetailRow) below row 3 so the
    # sheet's used range shrinks back down to row 3, without disturbing the
    # row-unbounded conditional formatting / data validation ranges already
    # defined on the sheet.
    if ($lastDetailRow -ge 4) {
        $ws.Range($ws.Cells.Item(4, 1), $ws.Cells.Item($lastDetailRow, 15)).EntireRow.ClearContents()
    }
}

# ---------------------------------------------------------------------------
# Helper: add a brand-new row 3 on sheets that currently only have the header
# rows (1 & 2), turning them into a single schema-reference row.
# ---------------------------------------------------------------------------
function Add-SchemaRow($ws, $section, $schemaName) {
    $ws.Cells.Item(3, 1).Value = $section
    $ws.Cells.Item(3, 2).Value = $schemaName
    $ws.Cells.Item(3, 5).Value = "schema"
    $ws.Cells.Item(3, 7).Value = $schemaName
    $ws.Cells.Item(3, 9).Value = "Yes"
}

# --- "Body" sheet: body/dateTime (+ body/commandRef) -> single schema row ---
$wsBody = $wb.Worksheets.Item("Body")
Set-SchemaRow $wsBody 4 "rejectParticipantOperation.210702Request"

# --- "200" sheet: content/dateTime (+ content/commandRef) -> schema row ---
$ws200 = $wb.Worksheets.Item("200")
Set-SchemaRow $ws200 4 "rejectParticipantOperation.210702Response"

# --- "204" sheet: empty -> single schema row referencing the Response schema ---
$ws204 = $wb.Worksheets.Item("204")
Add-SchemaRow $ws204 "content" "rejectParticipantOperation.210702Response"

# --- "400" sheet: content/dateTime + errorCode + errorCodeDescription + requestId -> schema row ---
$ws400 = $wb.Worksheets.Item("400")
Set-SchemaRow $ws400 6 "errorResponse"

# --- "401", "403", "404", "429", "500" sheets: empty -> schema row referencing errorResponse1 ---
foreach ($sheetName in @("401", "403", "404", "429", "500")) {
    $wsErr = $wb.Worksheets.Item($sheetName)
    Add-SchemaRow $wsErr "content" "errorResponse1"
}
